# Add a "Metadata" worksheet describing the columns of Sheet1, and make it
# the active sheet/tab, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- Add the new "Metadata" worksheet after Sheet1 ---
$meta = $wb.Worksheets.Add($null, $sheet1)
$meta.Name = "Metadata"

# --- Populate header row ---
$meta.Range("A1").Value = "Column Name"
$meta.Range("B1").Value = "Description"

# --- Populate column descriptions ---
$meta.Range("A2").Value = "pop_code"
$meta.Range("B2").Value = "Three-letter code for accession used"

$meta.Range("A3").Value = "rep"
$meta.Range("B3").Value = "Replicate, corresponding to a separate plant from which one leaf was measured"

$meta.Range("A4").Value = "leaf_area_cm2"
$meta.Range("B4").Value = "Leaf area in cm^2 from scans taken on a document scanner and manually measured in ImageJ. Leaves were used in leaf water drop adhesion assay, succulence measurement"

# --- Update Sheet1's selection (it's no longer the tab-selected sheet) ---
$sheet1.Range("A1:C1").Select()

# --- Select a cell on the Metadata sheet and make it the active tab ---
$meta.Range("C13").Select()
$meta.Activate()
$excel.ActiveWindow.Zoom = 160
